# Rename/restructure the eras_programs / non_eras_programs loop tags.
#
# Net effect (document order):
#   {#eras_programs}      -> {#e_progs}     (loop "e_progs" opens here, as before)
#   {#non_eras_programs}  -> {/e_progs}     (this old opening tag now CLOSES "e_progs")
#   {/eras_programs}      -> {#ne_progs}    (this old closing tag now OPENS "ne_progs")
#   {/non_eras_programs}  -> {/ne_progs}    (loop "ne_progs" closes here, as before)
#
# Each of the four literal tag strings occurs exactly once in the document, so a
# straightforward Find/Execute (case-sensitive, whole string, no wildcards) on each
# is unambiguous and safe regardless of order.

$d = $word.ActiveDocument

$wdReplaceOne = 1
$wdFindContinue = 1

$d.Content.Find.Execute(
    "#eras_programs}", $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false, "#e_progs}", $wdReplaceOne) | Out-Null

$d.Content.Find.Execute(
    "{#non_eras_programs}", $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false, "{/e_progs}", $wdReplaceOne) | Out-Null

$d.Content.Find.Execute(
    "{/eras_programs}", $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false, "{#ne_progs}", $wdReplaceOne) | Out-Null

$d.Content.Find.Execute(
    "{/non_eras_programs}", $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false, "{/ne_progs}", $wdReplaceOne) | Out-Null

Write-Output "done"
